$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '304.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.86%'
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '8'
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.18%'
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '8'
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.142'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.57%'
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '8'
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07589'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.96%'
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '8'
$ws.Range("G5").Style = "Normal"
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.331'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.14%'
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '8'
$ws.Range("G6").Style = "Normal"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.671'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.90%'
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '8'
$ws.Range("G7").Style = "Normal"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9413'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2.85%'
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '8'
$ws.Range("G8").Style = "Normal"
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1201'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.88%'
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '8'
$ws.Range("G9").Style = "Normal"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1823'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.21%'
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '8'
$ws.Range("G10").Style = "Normal"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09005'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.69%'
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '8'
$ws.Range("G11").Style = "Normal"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04154'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.37%'
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '8'
$ws.Range("G12").Style = "Normal"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1053'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.05%'
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '8'
$ws.Range("G13").Style = "Normal"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001293'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '2.29%'
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '8'
$ws.Range("G14").Style = "Normal"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005765'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.02%'
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '8'
$ws.Range("G15").Style = "Normal"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.341'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.17%'
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '8'
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.399'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.24%'
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '8'
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.3354'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.55%'
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '8'
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.615'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.71%'
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '8'
$ws.Range("G19").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-3.62%'
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '8'
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2807'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '3.34%'
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '8'
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.03859'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.75%'
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '8'
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001288'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.05%'
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '8'
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.003975'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-9.29%'
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '8'
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001354'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '5.63%'
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '8'
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003505'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-95.33%'
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '8'
$ws.Range("G26").Style = "Normal"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '8'
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '8'
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '8'
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '8'
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '8'
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '8'
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '8'
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '8'
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '8'
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '8'
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '8'
$ws.Range("G37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02401'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-3.35%'
$ws.Range("E38").Style = "Normal"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '8'
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05136'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-3.61%'
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '8'
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007719'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.55%'
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '8'
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1298'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.01%'
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '8'
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007556'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '16.17%'
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '8'
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003600'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '88.16%'
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '8'
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008222'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.26%'
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '8'
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3306'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.13%'
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '8'
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006847'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.99%'
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '8'
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.06%'
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '8'
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2723'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-24.92%'
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '8'
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004205'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '35.41%'
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '8'
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.06%'
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '8'
$ws.Range("G50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002002'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.06%'
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '8'
$ws.Range("G51").Style = "Normal"
